$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '265.49'
    'E2' = '1.68%'
    'D3' = '26.66'
    'E3' = '-1.33%'
    'D4' = '4.709'
    'E4' = '0.12%'
    'D5' = '0.06087'
    'E5' = '-1.53%'
    'D6' = '6.743'
    'E6' = '0.89%'
    'D7' = '0.8512'
    'E7' = '0.11%'
    'D8' = '0.9094'
    'E8' = '-0.58%'
    'D9' = '0.1407'
    'E9' = '-0.01%'
    'D10' = '0.05106'
    'E10' = '9.45%'
    'D11' = '0.07098'
    'E11' = '0.11%'
    'D12' = '0.03146'
    'E12' = '-0.13%'
    'D13' = '0.09022'
    'E13' = '-0.08%'
    'D14' = '0.001540'
    'E14' = '0.06%'
    'D15' = '0.0006042'
    'E15' = '-1.98%'
    'D16' = '0.005950'
    'E16' = '-0.91%'
    'D17' = '3.450'
    'E17' = '-0.01%'
    'E18' = '0.12%'
    'E19' = '0.78%'
    'E21' = '-1.40%'
    'D22' = '4.128'
    'E22' = '1.02%'
    'D23' = '0.04235'
    'E23' = '0.30%'
    'D24' = '0.001181'
    'E24' = '-2.93%'
    'D25' = '0.004059'
    'E25' = '6.75%'
    'E26' = '0.01%'
    'E27' = '6.57%'
    'D40' = '0.03926'
    'E40' = '1.05%'
    'D41' = '0.1115'
    'E41' = '0.28%'
    'D42' = '0.004198'
    'E42' = '2.54%'
    'D43' = '0.002106'
    'E43' = '-3.59%'
    'D44' = '0.01156'
    'E44' = '-29.17%'
    'D45' = '0.00005095'
    'E45' = '-1.24%'
    'D46' = '0.00000000750'
    'E46' = '0.01%'
    'D48' = '0.2587'
    'E48' = '55.20%'
    'D49' = '0.00002101'
    'E49' = '0.01%'
    'D50' = '0.0002001'
    'E50' = '0.01%'
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
